$d = $word.ActiveDocument
$newText = "Perioadele campaniei din Cygnus: 10-19 august, 9-18 septembrie, 8-17 octombrie"

# Collect the paragraphs that need to change first (indices shift as we edit,
# so resolve them by index from the back of the document forward).
$targets = @()
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Perioadele campaniei din 2018 pentru*") {
        $targets += $idx
    }
}

for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $p = $d.Paragraphs($targets[$i])
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Delete()
    $r.InsertAfter($newText)
}
